# "Updated car model drag coeff and added low drag configs"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")
$ws.Activate()

# Lift Coefficient CL: -2 -> -1.98
$ws.Range("C8").Value = -1.98

# Drag Coefficient CD: -1.2 -> -1.33 (low drag config)
$ws.Range("C9").Value = -1.33

# Front Aero Distribution: 47 -> computed low-drag value via formula
$ws.Range("C12").Formula = "=100-56.3"

# CL Scale Multiplier: 1.1 -> 1.15
$ws.Range("C13").Value = 1.1499999999999999

# Leave the cursor where the author left it
$ws.Range("E8").Select()
